$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "SBO_DEF" in column F, matching the style of the other headers
$ws.Range("F1").Value = "SBO_DEF"
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Fill F2:F14 with the literal value "[]"
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}
